$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-11 01:43:26"

for ($r = 2; $r -le 22; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
